# Progress - BGZ to FHIR.xlsx update
# Applies: new ZIB issue rows, resolution notes, mapping resolutions,
# and refreshed view/selection state (per commit "Updated mappings and progress .xlsx")

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Overview"
$ws2 = $wb.Worksheets.Item(2)   # "Issues"

# ---------------------------------------------------------------
# 1. Issues sheet (sheet2): add resolution remarks to existing rows
# ---------------------------------------------------------------
$ws2.Cells.Item(25, 6).Value = "Cardinaliteit gewijzigd van 1 naar 0..1"
$ws2.Cells.Item(26, 6).Value = "Cardinaliteit Organisatienaam gewijzigd van 1 naar 0..1"

# ---------------------------------------------------------------
# 2. Overview sheet (sheet1): fill in newly-resolved "Resolution" notes
# ---------------------------------------------------------------
$ws1.Cells.Item(20, 9).Value = "WebEx 15-05"
$ws1.Cells.Item(21, 9).Value = "HL7 WGM 31-03 / WebEx 15-05"
$ws1.Cells.Item(22, 9).Value = "WebEx 15-05"

# ---------------------------------------------------------------
# 3. Issues sheet (sheet2): append four new issue rows (35-38)
# ---------------------------------------------------------------
$ws2.Cells.Item(35, 2).Value = "ZIB-605 "
$ws2.Cells.Item(35, 3).Value = "Naam MedischHulpmiddel veranderen in bijvoorbeeld MedischHulpmiddelGebruik"
$ws2.Cells.Item(35, 1).Value = "MedicalAid"
$ws2.Cells.Item(35, 4).Value = "In intake"
$ws2.Cells.Item(35, 5).Value = "Unresolved"

$ws2.Cells.Item(36, 3).Value = "Naam wijzigen van Woonsituatie naar Woonvoorziening"
$ws2.Cells.Item(36, 2).Value = "ZIB-604"
$ws2.Cells.Item(36, 1).Value = "LivingSituatioin"
$ws2.Cells.Item(36, 4).Value = "In intake"
$ws2.Cells.Item(36, 5).Value = "Unresolved"

$ws2.Cells.Item(37, 2).Value = "#12685"
$ws2.Cells.Item(37, 3).Value = "Linking of Accounts required"
$ws2.Cells.Item(37, 1).Value = "Payor"
$ws2.Cells.Item(37, 5).Value = "This was discussed and in concept agree, but we would like more detail on the usage, potential types, and impact on the linked to account. (e.g. Does the balance change?)"
$ws2.Cells.Item(37, 4).Value = "Triaged"

$ws2.Cells.Item(38, 2).Value = "#13415"
$ws2.Cells.Item(38, 1).Value = "VoedingAdvies"
$ws2.Cells.Item(38, 3).Value = "Merge NutritionOrder texture and fluidConsistencyType "
$ws2.Cells.Item(38, 4).Value = "Submitted"

# ---------------------------------------------------------------
# 4. Issues sheet (sheet2): column A is a touch narrower now that
#    longer entries no longer need "best fit" sizing
# ---------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 23.75

# ---------------------------------------------------------------
# 5. Refresh window state: Overview tab is now the active / selected
#    tab, with an updated zoom level and selection; Issues keeps a
#    new scroll position, zoom level and selection as well.
# ---------------------------------------------------------------
$ws2.Activate()
$w2 = $excel.ActiveWindow
$w2.Zoom = 130
$ws2.Range("B38:C38").Select()

$ws1.Activate()
$w1 = $excel.ActiveWindow
$w1.Zoom = 100
$ws1.Range("F22").Select()
